$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swerve Home table (columns E:G, rows 1-5) ---
# Column E (labels) filled top-to-bottom first
$ws.Range("E1").Value = "Swerve Home"
$ws.Range("E2").Value = "Front left"
$ws.Range("E3").Value = "Front right"
$ws.Range("E4").Value = "Back left"
$ws.Range("E5").Value = "Back right"

# Then the header row for columns F and G
$ws.Range("F1").Value = "Home Relative to Halsensor"
$ws.Range("G1").Value = "Halsensor position relative to home"

# Then the numeric data + formula for each row
$ws.Range("F2").Value = -8.9524135589599592
$ws.Range("G2").Formula = "=-F2+9"

$ws.Range("F3").Value = -8.7857446670532209
$ws.Range("G3").Formula = "=-F3+9"

$ws.Range("F4").Value = -8.9047937393188406
$ws.Range("G4").Formula = "=-F4+9"

$ws.Range("F5").Value = -8.9524135589599592
$ws.Range("G5").Formula = "=-F5+9"

# --- Swerve drive PID block (columns E:F, rows 8-12) ---
$ws.Range("E8").Value = "Swerve drive PID"

$ws.Range("E9").Value = "P"
$ws.Range("F9").Value = 0.00001

$ws.Range("E10").Value = "I"

$ws.Range("E11").Value = "D"

$ws.Range("E12").Value = "F"
$ws.Range("F12").Value = 0.000166

# --- Column widths for the new columns (best-fit, like Excel's AutoFit) ---
# Values chosen so the engine's pixel-quantized stored width lands on the
# closest representable value to the original bestFit widths
# (12.33203125 / 23.9296875 / 30.86328125 character units).
$ws.Columns.Item(5).ColumnWidth = 11.5
$ws.Columns.Item(6).ColumnWidth = 23.166666666666668
$ws.Columns.Item(7).ColumnWidth = 30

# --- Restore active selection to match the final saved state ---
$ws.Range("H14").Select() | Out-Null
